# prob28 (p&s part) insertion
# - Inserts a fresh 10-row gap before row 141 (old rows 141-171 slide down to 151-181)
# - Inserts a second 10-row gap before the (now) row 171 (old row 161 -> 181, old row 171 -> 191)
# - Fills the newly opened rows 165-168 with the x0015..x0018 key/value pairs for prob28
# - Restores the sheet selection to match the author's final cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Make room for the new problem block right after the existing x0001-x0014 block.
$ws.Rows("141:150").Insert()
# Make room between the relocated y0001/z0001 marker rows so they land on 181/191.
$ws.Rows("171:180").Insert()

# New key/value rows for prob28 (x0015-x0018).
# Column B is written before column A on row 167 so the shared-string table
# picks up entries in the same order as the source workbook.
$ws.Range("A165").Value = "x0015"
$ws.Range("B165").Value = "함수의 정의역의 각 원소에 대한 함숫값의 범위를 조사합니다."

$ws.Range("A166").Value = "x0016"
$ws.Range("B166").Value = "조건에 맞는 가능한 치역을 모두 구합니다. "

$ws.Range("B167").Value = "각 치역의 경우에 대해서 조건을 만족시키는 함수의 개수를 구합니다."
$ws.Range("A167").Value = "x0017"

$ws.Range("A168").Value = "x0018"
$ws.Range("B168").Value = "각각의 개수를 모두 더해서 조건을 만족시키는 전체 개수를 구합니다. "

# Match the author's final view/selection state.
$ws.Activate()
$ws.Range("B176").Select()
$excel.ActiveWindow.ScrollRow = 151
$excel.ActiveWindow.ScrollColumn = 1
